$d = $word.ActiveDocument

# 1. The label run "URL to Your Coding Assignment Video:" gets a trailing
#    space appended. Using Find/Execute lets Word manage the
#    xml:space="preserve" attribute for us automatically.
$d.Content.Find.Execute("URL to Your Coding Assignment Video:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "URL to Your Coding Assignment Video: ", 2)

# 2. Re-locate that paragraph (its text now carries the trailing space) and
#    append a second run holding "???" right before the paragraph mark.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "URL to Your Coding Assignment Video:*") {
        $targetIndex = $i
        break
    }
}

$paraRange = $d.Paragraphs($targetIndex).Range
$paraRange.InsertAfter("???")

# 3. The newly typed text inherited the bold formatting of the run it was
#    appended to (standard Word behavior). Turn bold back off for just the
#    "???" text so it becomes its own, non-bold run - matching the
#    GitHub-URL paragraph just above it, which uses the same bold-label /
#    non-bold-value pattern.
$paraRangeAfter = $d.Paragraphs($targetIndex).Range
$newRunRange = $d.Range($paraRangeAfter.End - 4, $paraRangeAfter.End - 1)
$newRunRange.Font.Bold = $false
$newRunRange.Font.BoldBi = $false
